# Add the All-Staff meeting dates (column D) and reminder dates (column E)
# to the "Facilitator Tracker" sheet. Meetings are the 3rd Monday of each
# month at 1:00-2:30 PM, except February which moved to 2/17 (Tue) because
# the 3rd Monday (2/16) is a holiday. Reminder dates are 2 weeks prior.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Facilitator Tracker")

# row -> (All-Staff Date, Reminder Date)
$dates = @(
    @{ Row = 2;  AllStaff = "1/19";       Reminder = "1/5"  },  # January 2026
    @{ Row = 3;  AllStaff = "2/17 (Tue)"; Reminder = "2/3"  },  # February 2026 (moved from 2/16 holiday)
    @{ Row = 4;  AllStaff = "3/16";       Reminder = "3/2"  },  # March 2026
    @{ Row = 5;  AllStaff = "4/20";       Reminder = "4/6"  },  # April 2026
    @{ Row = 6;  AllStaff = "5/18";       Reminder = "5/4"  },  # May 2026
    @{ Row = 7;  AllStaff = "6/15";       Reminder = "6/1"  },  # June 2026
    @{ Row = 8;  AllStaff = "7/20";       Reminder = "7/6"  },  # July 2026
    @{ Row = 9;  AllStaff = "8/17";       Reminder = "8/3"  },  # August 2026
    @{ Row = 10; AllStaff = "9/21";       Reminder = "9/7"  },  # September 2026
    @{ Row = 11; AllStaff = "10/19";      Reminder = "10/5" },  # October 2026
    @{ Row = 12; AllStaff = "11/16";      Reminder = "11/2" },  # November 2026
    @{ Row = 13; AllStaff = "12/21";      Reminder = "12/7" }   # December 2026
)

foreach ($entry in $dates) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.AllStaff   # column D: All-Staff Date
    $ws.Cells.Item($entry.Row, 5).Value = $entry.Reminder   # column E: Reminder Date
}
